$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7829526662826538
$ws.Range("B1").Value = 2.165490388870239
$ws.Range("D1").Value = 1.08585786819458
$ws.Range("E1").Value = 0.5727079510688782
